$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1143.4615
$ws.Range("J17").Value = 1154
$ws.Range("L17").Value = 3462
$ws.Range("N17").Value = -3798
$ws.Range("H28").Value = 5907.1113
$ws.Range("I28").Value = 6226.1177
$ws.Range("J28").Value = 484
$ws.Range("K28").Value = 6226.1177
$ws.Range("L28").Value = 484
$ws.Range("M28").Value = -5741.1177
$ws.Range("N28").Value = -1454
$ws.Range("H32").Value = 2849.5
$ws.Range("I32").Value = 1699
$ws.Range("K32").Value = 1699
$ws.Range("M32").Value = -1373
$ws.Range("H33").Value = 2117.6
$ws.Range("I33").Value = 219.2
$ws.Range("K33").Value = 219.2
$ws.Range("M33").Value = 9.800000000000011
$ws.Range("H43").Value = 4363.8076
$ws.Range("I43").Value = 3756.1428
$ws.Range("J43").Value = 4587.684
$ws.Range("K43").Value = 3756.1428
$ws.Range("L43").Value = 4587.684
$ws.Range("M43").Value = -3687.1428
$ws.Range("N43").Value = -4725.684
$ws.Range("H62").Value = 50518.61
$ws.Range("I62").Value = 76307.73
$ws.Range("J62").Value = 9992.857
$ws.Range("K62").Value = 76307.73
$ws.Range("L62").Value = 9992.857
$ws.Range("M62").Value = -75683.73
$ws.Range("N62").Value = -11240.857
$ws.Range("H65").Value = 50518.61
$ws.Range("I65").Value = 76307.73
$ws.Range("J65").Value = 9992.857
$ws.Range("K65").Value = 381538.65
$ws.Range("L65").Value = 49964.285
$ws.Range("M65").Value = -378418.65
$ws.Range("N65").Value = -56204.285
$ws.Range("H70").Value = 10007550
$ws.Range("J70").Value = 8389.223
$ws.Range("L70").Value = 25167.669
$ws.Range("N70").Value = -25707.669
$ws.Range("H73").Value = 10007550
$ws.Range("J73").Value = 8389.223
$ws.Range("L73").Value = 25167.669
$ws.Range("N73").Value = -27039.669
$ws.Range("H75").Value = 105000
$ws.Range("J75").Value = 105000
$ws.Range("L75").Value = 105000
$ws.Range("N75").Value = -106872
$ws.Range("H78").Value = 105000
$ws.Range("J78").Value = 105000
$ws.Range("L78").Value = 315000
$ws.Range("N78").Value = -324360
$ws.Range("H98").Value = 7392.875
$ws.Range("I98").Value = 10072.823
$ws.Range("J98").Value = 884.4286
$ws.Range("K98").Value = 10072.823
$ws.Range("L98").Value = 884.4286
$ws.Range("M98").Value = -8574.823
$ws.Range("N98").Value = -3880.4286
$ws.Range("H100").Value = 5824.875
$ws.Range("I100").Value = 2742.8572
$ws.Range("J100").Value = 8222
$ws.Range("K100").Value = 2742.8572
$ws.Range("L100").Value = 8222
$ws.Range("M100").Value = -2201.8572
$ws.Range("N100").Value = -9304
$ws.Range("H122").Value = 7392.875
$ws.Range("I122").Value = 10072.823
$ws.Range("J122").Value = 884.4286
$ws.Range("K122").Value = 30218.469
$ws.Range("L122").Value = 2653.2858
$ws.Range("M122").Value = -27768.469
$ws.Range("N122").Value = -7553.2858
$ws.Range("H123").Value = 50000
$ws.Range("I123").Value = 30000
$ws.Range("J123").Value = 70000
$ws.Range("K123").Value = 30000
$ws.Range("L123").Value = 70000
$ws.Range("M123").Value = -25100
$ws.Range("N123").Value = -79800
$ws.Range("H132").Value = 4015.75
$ws.Range("I132").Value = 4242
$ws.Range("J132").Value = 3699
$ws.Range("K132").Value = 12726
$ws.Range("L132").Value = 11097
$ws.Range("M132").Value = -10196
$ws.Range("N132").Value = -16157
$ws.Range("H135").Value = 1177.9445
$ws.Range("I135").Value = 1170.9286
$ws.Range("J135").Value = 1202.5
$ws.Range("K135").Value = 10538.3574
$ws.Range("L135").Value = 10822.5
$ws.Range("M135").Value = -8003.357399999999
$ws.Range("N135").Value = -15892.5
$ws.Range("H137").Value = 1320.8572
$ws.Range("I137").Value = 1311.8064
$ws.Range("J137").Value = 1391
$ws.Range("K137").Value = 3935.4192
$ws.Range("L137").Value = 4173
$ws.Range("M137").Value = -1385.4192
$ws.Range("N137").Value = -9273
$ws.Range("H140").Value = 60000
$ws.Range("I140").Value = 40000
$ws.Range("J140").Value = 70000
$ws.Range("K140").Value = 40000
$ws.Range("L140").Value = 70000
$ws.Range("M140").Value = -34820
$ws.Range("N140").Value = -80360
$ws.Range("H141").Value = 4424.4
$ws.Range("J141").Value = 3874.25
$ws.Range("L141").Value = 11622.75
$ws.Range("N141").Value = -21982.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2599.889
$ws.Range("I2").Value = 1618
$ws.Range("K2").Value = 1618
$ws.Range("M2").Value = -1505
$ws.Range("H32").Value = 3892.6035
$ws.Range("I32").Value = 3017.7454
$ws.Range("J32").Value = 19931.666
$ws.Range("K32").Value = 3017.7454
$ws.Range("L32").Value = 19931.666
$ws.Range("M32").Value = -2730.7454
$ws.Range("N32").Value = -20505.666
$ws.Range("H33").Value = 120000
$ws.Range("I33").Value = 120000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 120000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -119671
$ws.Range("N33").ClearContents()
$ws.Range("H36").Value = 9526
$ws.Range("I36").Value = 9526
$ws.Range("K36").Value = 9526
$ws.Range("M36").Value = -9180
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H45").Value = 2004
$ws.Range("I45").Value = 1802
$ws.Range("J45").Value = 3014
$ws.Range("K45").Value = 1802
$ws.Range("L45").Value = 3014
$ws.Range("M45").Value = -1425
$ws.Range("N45").Value = -3768
$ws.Range("H61").Value = 3568.75
$ws.Range("I61").Value = 3568.75
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3568.75
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3356.75
$ws.Range("N61").ClearContents()
$ws.Range("H74").Value = 1847.9131
$ws.Range("I74").Value = 1847.0476
$ws.Range("K74").Value = 1847.0476
$ws.Range("M74").Value = -973.0476000000001
$ws.Range("H77").Value = 1847.9131
$ws.Range("I77").Value = 1847.0476
$ws.Range("K77").Value = 9235.238000000001
$ws.Range("M77").Value = -4867.238000000001
$ws.Range("H116").Value = 2599.889
$ws.Range("I116").Value = 1618
$ws.Range("K116").Value = 1618
$ws.Range("M116").Value = 676
$ws.Range("H122").Value = 1971.5
$ws.Range("I122").Value = 1636.7693
$ws.Range("K122").Value = 4910.3079
$ws.Range("M122").Value = -2460.3079
$ws.Range("H132").Value = 2846.6667
$ws.Range("I132").Value = 2839.5
$ws.Range("J132").Value = 2990
$ws.Range("K132").Value = 8518.5
$ws.Range("L132").Value = 8970
$ws.Range("M132").Value = -5988.5
$ws.Range("N132").Value = -14030
$ws.Range("H136").Value = 3568.75
$ws.Range("I136").Value = 3568.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10706.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8156.25
$ws.Range("N136").ClearContents()
$ws.Range("H137").Value = 104000
$ws.Range("J137").Value = 104000
$ws.Range("L137").Value = 104000
$ws.Range("N137").Value = -114200

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H2").Value = 54196.332
$ws.Range("J2").Value = 54196.332
$ws.Range("L2").Value = 54196.332
$ws.Range("N2").Value = -54422.332
$ws.Range("H3").Value = 2599.889
$ws.Range("I3").Value = 1618
$ws.Range("K3").Value = 1618
$ws.Range("M3").Value = -1504
$ws.Range("H22").Value = 957.6
$ws.Range("I22").Value = 1129.4286
$ws.Range("J22").Value = 556.6667
$ws.Range("K22").Value = 1129.4286
$ws.Range("L22").Value = 556.6667
$ws.Range("M22").Value = -956.4286
$ws.Range("N22").Value = -902.6667
$ws.Range("H50").Value = 99950
$ws.Range("J50").Value = 99950
$ws.Range("L50").Value = 99950
$ws.Range("N50").Value = -101098
$ws.Range("H62").Value = 40000
$ws.Range("I62").Value = 30000
$ws.Range("J62").Value = 45000
$ws.Range("K62").Value = 30000
$ws.Range("L62").Value = 45000
$ws.Range("M62").Value = -29314
$ws.Range("N62").Value = -46372
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H65").Value = 40000
$ws.Range("I65").Value = 30000
$ws.Range("J65").Value = 45000
$ws.Range("K65").Value = 90000
$ws.Range("L65").Value = 135000
$ws.Range("M65").Value = -86568
$ws.Range("N65").Value = -141864
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H94").Value = 3418
$ws.Range("I94").Value = 2759.8
$ws.Range("K94").Value = 2759.8
$ws.Range("M94").Value = -2308.8
$ws.Range("H96").Value = 6874.7
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()
$ws.Range("H134").Value = 2705.3333
$ws.Range("I134").Value = 2668.5
$ws.Range("K134").Value = 8005.5
$ws.Range("M134").Value = -5470.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 7799.8
$ws.Range("I16").Value = 6000
$ws.Range("J16").Value = 10499.5
$ws.Range("K16").Value = 6000
$ws.Range("L16").Value = 10499.5
$ws.Range("M16").Value = -5713
$ws.Range("N16").Value = -11073.5
$ws.Range("H19").Value = 2501
$ws.Range("I19").Value = 403
$ws.Range("J19").Value = 3550
$ws.Range("K19").Value = 403
$ws.Range("L19").Value = 3550
$ws.Range("M19").Value = -233
$ws.Range("N19").Value = -3890
$ws.Range("H24").Value = 2501
$ws.Range("I24").Value = 403
$ws.Range("J24").Value = 3550
$ws.Range("K24").Value = 403
$ws.Range("L24").Value = 3550
$ws.Range("M24").Value = -233
$ws.Range("N24").Value = -3890
$ws.Range("H31").Value = 5520.5
$ws.Range("I31").Value = 1921.4445
$ws.Range("J31").Value = 11998.8
$ws.Range("K31").Value = 1921.4445
$ws.Range("L31").Value = 11998.8
$ws.Range("M31").Value = -1626.4445
$ws.Range("N31").Value = -12588.8
$ws.Range("H34").Value = 5520.5
$ws.Range("I34").Value = 1921.4445
$ws.Range("J34").Value = 11998.8
$ws.Range("K34").Value = 1921.4445
$ws.Range("L34").Value = 11998.8
$ws.Range("M34").Value = -1719.4445
$ws.Range("N34").Value = -12402.8
$ws.Range("H41").Value = 36459.08
$ws.Range("I41").Value = 6264.75
$ws.Range("J41").Value = 49878.777
$ws.Range("K41").Value = 6264.75
$ws.Range("L41").Value = 49878.777
$ws.Range("M41").Value = -5836.75
$ws.Range("N41").Value = -50734.777
$ws.Range("H55").Value = 120000
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H56").Value = 24495
$ws.Range("I56").Value = 24495
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 24495
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -23650
$ws.Range("N56").ClearContents()
$ws.Range("H58").Value = 1123.909
$ws.Range("I58").Value = 1136.3
$ws.Range("J58").Value = 1000
$ws.Range("K58").Value = 1136.3
$ws.Range("L58").Value = 1000
$ws.Range("M58").Value = -933.3
$ws.Range("N58").Value = -1406
$ws.Range("H59").Value = 59612.125
$ws.Range("I59").Value = 35000
$ws.Range("J59").Value = 63128.145
$ws.Range("K59").Value = 35000
$ws.Range("L59").Value = 63128.145
$ws.Range("M59").Value = -33855
$ws.Range("N59").Value = -65418.145
$ws.Range("H60").Value = 40949.1
$ws.Range("I60").Value = 29833
$ws.Range("J60").Value = 45713.145
$ws.Range("K60").Value = 29833
$ws.Range("L60").Value = 45713.145
$ws.Range("M60").Value = -29322
$ws.Range("N60").Value = -46735.145
$ws.Range("H62").Value = 14844.286
$ws.Range("J62").Value = 16401
$ws.Range("L62").Value = 16401
$ws.Range("N62").Value = -17649
$ws.Range("H65").Value = 14844.286
$ws.Range("J65").Value = 16401
$ws.Range("L65").Value = 82005
$ws.Range("N65").Value = -88245
$ws.Range("H99").Value = 1431077
$ws.Range("I99").Value = 2002788
$ws.Range("J99").Value = 1799.5
$ws.Range("K99").Value = 2002788
$ws.Range("L99").Value = 1799.5
$ws.Range("M99").Value = -2001290
$ws.Range("N99").Value = -4795.5
$ws.Range("H105").Value = 2867.818
$ws.Range("I105").Value = 2456
$ws.Range("J105").Value = 3966
$ws.Range("K105").Value = 2456
$ws.Range("L105").Value = 3966
$ws.Range("M105").Value = -709
$ws.Range("N105").Value = -7460
$ws.Range("H113").Value = 7799.8
$ws.Range("I113").Value = 6000
$ws.Range("J113").Value = 10499.5
$ws.Range("K113").Value = 6000
$ws.Range("L113").Value = 10499.5
$ws.Range("M113").Value = -3830
$ws.Range("N113").Value = -14839.5
$ws.Range("H126").Value = 1431077
$ws.Range("I126").Value = 2002788
$ws.Range("J126").Value = 1799.5
$ws.Range("K126").Value = 6008364
$ws.Range("L126").Value = 5398.5
$ws.Range("M126").Value = -6005894
$ws.Range("N126").Value = -10338.5
$ws.Range("H132").Value = 1947.4445
$ws.Range("I132").Value = 1646.8572
$ws.Range("K132").Value = 4940.571599999999
$ws.Range("M132").Value = -2410.571599999999
$ws.Range("H134").Value = 3199.182
$ws.Range("I134").Value = 2644.15
$ws.Range("K134").Value = 7932.450000000001
$ws.Range("M134").Value = -5397.450000000001
$ws.Range("H136").Value = 1123.909
$ws.Range("I136").Value = 1136.3
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 3408.9
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -858.8999999999996
$ws.Range("N136").Value = -8100
$ws.Range("H138").Value = 367500
$ws.Range("J138").Value = 367500
$ws.Range("L138").Value = 367500
$ws.Range("N138").Value = -377780

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100.84615
$ws.Range("I2").Value = 89.125
$ws.Range("J2").Value = 119.6
$ws.Range("K2").Value = 534.75
$ws.Range("L2").Value = 717.5999999999999
$ws.Range("M2").Value = -421.75
$ws.Range("N2").Value = -943.5999999999999
$ws.Range("H8").Value = 213
$ws.Range("I8").Value = 213
$ws.Range("K8").Value = 639
$ws.Range("M8").Value = -500
$ws.Range("H9").Value = 253666.33
$ws.Range("J9").Value = 253666.33
$ws.Range("L9").Value = 760998.99
$ws.Range("N9").Value = -761446.99
$ws.Range("H34").Value = 540
$ws.Range("J34").Value = 562.5
$ws.Range("L34").Value = 1687.5
$ws.Range("N34").Value = -1855.5
$ws.Range("H81").Value = 3733.2144
$ws.Range("I81").Value = 708.6667
$ws.Range("J81").Value = 6001.625
$ws.Range("K81").Value = 2126.0001
$ws.Range("L81").Value = 18004.875
$ws.Range("M81").Value = -1003.0001
$ws.Range("N81").Value = -20250.875
$ws.Range("H84").Value = 3733.2144
$ws.Range("I84").Value = 708.6667
$ws.Range("J84").Value = 6001.625
$ws.Range("K84").Value = 6378.0003
$ws.Range("L84").Value = 54014.625
$ws.Range("M84").Value = -762.0002999999997
$ws.Range("N84").Value = -65246.625
$ws.Range("H97").Value = 773.4286
$ws.Range("J97").Value = 934.75
$ws.Range("L97").Value = 2804.25
$ws.Range("N97").Value = -3796.25
$ws.Range("H100").Value = 35862.5
$ws.Range("I100").Value = 5025
$ws.Range("K100").Value = 15075
$ws.Range("M100").Value = -14264
$ws.Range("H110").Value = 7856.4287
$ws.Range("I110").Value = 4999.6665
$ws.Range("J110").Value = 9999
$ws.Range("K110").Value = 14998.9995
$ws.Range("L110").Value = 29997
$ws.Range("M110").Value = -10908.9995
$ws.Range("N110").Value = -38177
$ws.Range("H112").Value = 47620428
$ws.Range("I112").Value = 166667250
$ws.Range("J112").Value = 1701
$ws.Range("K112").Value = 500001750
$ws.Range("L112").Value = 5103
$ws.Range("M112").Value = -500000642
$ws.Range("N112").Value = -7319
$ws.Range("H117").Value = 1692.6
$ws.Range("J117").Value = 1328
$ws.Range("L117").Value = 3984
$ws.Range("N117").Value = -10868
$ws.Range("H121").Value = 757.875
$ws.Range("I121").Value = 317.375
$ws.Range("K121").Value = 952.125
$ws.Range("M121").Value = 357.875
$ws.Range("H122").Value = 2982.875
$ws.Range("J122").Value = 2982.875
$ws.Range("L122").Value = 26845.875
$ws.Range("N122").Value = -31745.875
$ws.Range("H124").Value = 1900
$ws.Range("I124").Value = 1900
$ws.Range("K124").Value = 5700
$ws.Range("M124").Value = -790
$ws.Range("H129").Value = 2951.7693
$ws.Range("J129").Value = 4575
$ws.Range("L129").Value = 13725
$ws.Range("N129").Value = -23725
$ws.Range("H131").Value = 2529.1667
$ws.Range("I131").Value = 1475
$ws.Range("J131").Value = 2740
$ws.Range("K131").Value = 4425
$ws.Range("L131").Value = 8220
$ws.Range("M131").Value = 615
$ws.Range("N131").Value = -18300
$ws.Range("H132").Value = 2719.9
$ws.Range("I132").Value = 2439.8
$ws.Range("K132").Value = 21958.2
$ws.Range("M132").Value = -19428.2

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 22500
$ws.Range("I26").Value = 10000
$ws.Range("J26").Value = 35000
$ws.Range("K26").Value = 10000
$ws.Range("L26").Value = 35000
$ws.Range("M26").Value = -9720
$ws.Range("N26").Value = -35560
$ws.Range("H50").Value = 22500
$ws.Range("I50").Value = 10000
$ws.Range("J50").Value = 35000
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 35000
$ws.Range("M50").Value = -9502
$ws.Range("N50").Value = -35996
$ws.Range("H93").Value = 26049.4
$ws.Range("J93").Value = 26049.4
$ws.Range("L93").Value = 26049.4
$ws.Range("N93").Value = -29793.4
$ws.Range("H102").Value = 13198.8
$ws.Range("I102").Value = 8665
$ws.Range("J102").Value = 19999.5
$ws.Range("K102").Value = 8665
$ws.Range("L102").Value = 19999.5
$ws.Range("M102").Value = -7043
$ws.Range("N102").Value = -23243.5
$ws.Range("H113").Value = 27783372
$ws.Range("J113").Value = 9699.866
$ws.Range("L113").Value = 9699.866
$ws.Range("N113").Value = -14039.866
$ws.Range("H122").Value = 35791.25
$ws.Range("I122").Value = 39944
$ws.Range("J122").Value = 23333
$ws.Range("K122").Value = 119832
$ws.Range("L122").Value = 69999
$ws.Range("M122").Value = -117382
$ws.Range("N122").Value = -74899
$ws.Range("H126").Value = 3050
$ws.Range("I126").Value = 3050
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9150
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6680
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2997.3333
$ws.Range("I132").Value = 2944.4736
$ws.Range("K132").Value = 8833.4208
$ws.Range("M132").Value = -6303.4208
$ws.Range("H134").Value = 42074.11
$ws.Range("J134").Value = 42074.11
$ws.Range("L134").Value = 126222.33
$ws.Range("N134").Value = -131292.33

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7112.4287
$ws.Range("I7").Value = 5849.5
$ws.Range("K7").Value = 5849.5
$ws.Range("M7").Value = -5737.5
$ws.Range("H9").Value = 324
$ws.Range("I9").Value = 68
$ws.Range("J9").Value = 580
$ws.Range("K9").Value = 68
$ws.Range("L9").Value = 580
$ws.Range("M9").Value = 156
$ws.Range("N9").Value = -1028
$ws.Range("H16").Value = 1377.4
$ws.Range("J16").Value = 2170.75
$ws.Range("L16").Value = 2170.75
$ws.Range("N16").Value = -2510.75
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H30").Value = 16
$ws.Range("I30").Value = 16
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 16
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = 92
$ws.Range("N30").ClearContents()
$ws.Range("H40").Value = 4952.625
$ws.Range("I40").Value = 3417.5334
$ws.Range("K40").Value = 3417.5334
$ws.Range("M40").Value = -3281.5334
$ws.Range("H43").Value = 13477.378
$ws.Range("J43").Value = 14280.518
$ws.Range("L43").Value = 14280.518
$ws.Range("N43").Value = -14666.518
$ws.Range("H46").Value = 1543.4348
$ws.Range("I46").Value = 1121.2354
$ws.Range("J46").Value = 1790.931
$ws.Range("K46").Value = 1121.2354
$ws.Range("L46").Value = 1790.931
$ws.Range("M46").Value = -933.2354
$ws.Range("N46").Value = -2166.931
$ws.Range("H55").Value = 1564.4706
$ws.Range("I55").Value = 236.54546
$ws.Range("K55").Value = 236.54546
$ws.Range("M55").Value = -63.54545999999999
$ws.Range("H68").Value = 5414.8335
$ws.Range("J68").Value = 7903.75
$ws.Range("L68").Value = 7903.75
$ws.Range("N68").Value = -9401.75
$ws.Range("H71").Value = 5414.8335
$ws.Range("J71").Value = 7903.75
$ws.Range("L71").Value = 39518.75
$ws.Range("N71").Value = -47006.75
$ws.Range("H82").Value = 2962.8235
$ws.Range("I82").Value = 1518.7778
$ws.Range("J82").Value = 4587.375
$ws.Range("K82").Value = 1518.7778
$ws.Range("L82").Value = 4587.375
$ws.Range("M82").Value = -1157.7778
$ws.Range("N82").Value = -5309.375
$ws.Range("H85").Value = 2962.8235
$ws.Range("I85").Value = 1518.7778
$ws.Range("J85").Value = 4587.375
$ws.Range("K85").Value = 1518.7778
$ws.Range("L85").Value = 4587.375
$ws.Range("M85").Value = -270.7778000000001
$ws.Range("N85").Value = -7083.375
$ws.Range("H92").Value = 47694.5
$ws.Range("J92").Value = 47694.5
$ws.Range("L92").Value = 47694.5
$ws.Range("N92").Value = -52686.5
$ws.Range("H93").Value = 5375.25
$ws.Range("I93").Value = 2286.2856
$ws.Range("K93").Value = 2286.2856
$ws.Range("M93").Value = -1038.2856
$ws.Range("H100").Value = 5217.909
$ws.Range("I100").Value = 3119.6
$ws.Range("J100").Value = 9714.286
$ws.Range("K100").Value = 3119.6
$ws.Range("L100").Value = 9714.286
$ws.Range("M100").Value = -2578.6
$ws.Range("N100").Value = -10796.286
$ws.Range("H122").Value = 4164.0625
$ws.Range("I122").Value = 5180.6665
$ws.Range("J122").Value = 2857
$ws.Range("K122").Value = 15541.9995
$ws.Range("L122").Value = 8571
$ws.Range("M122").Value = -13091.9995
$ws.Range("N122").Value = -13471
$ws.Range("H126").Value = 7112.4287
$ws.Range("I126").Value = 5849.5
$ws.Range("K126").Value = 17548.5
$ws.Range("M126").Value = -15078.5
$ws.Range("H132").Value = 3827.3572
$ws.Range("I132").Value = 3659.8635
$ws.Range("J132").Value = 4441.5
$ws.Range("K132").Value = 10979.5905
$ws.Range("L132").Value = 13324.5
$ws.Range("M132").Value = -8449.5905
$ws.Range("N132").Value = -18384.5
$ws.Range("H135").Value = 300000
$ws.Range("J135").Value = 300000
$ws.Range("L135").Value = 300000
$ws.Range("N135").Value = -310140
$ws.Range("H136").Value = 12256.135
$ws.Range("I136").Value = 1781.8334
$ws.Range("J136").Value = 15398.425
$ws.Range("K136").Value = 5345.5002
$ws.Range("L136").Value = 46195.27499999999
$ws.Range("M136").Value = -2795.5002
$ws.Range("N136").Value = -51295.27499999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 6903334
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H81").Value = 2684.1538
$ws.Range("I81").Value = 1554.8572
$ws.Range("J81").Value = 4001.6667
$ws.Range("K81").Value = 3109.7144
$ws.Range("L81").Value = 8003.3334
$ws.Range("M81").Value = -2048.7144
$ws.Range("N81").Value = -10125.3334
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H84").Value = 2684.1538
$ws.Range("I84").Value = 1554.8572
$ws.Range("J84").Value = 4001.6667
$ws.Range("K84").Value = 15548.572
$ws.Range("L84").Value = 40016.667
$ws.Range("M84").Value = -10244.572
$ws.Range("N84").Value = -50624.667
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H113").Value = 1530.75
$ws.Range("I113").Value = 1239.4445
$ws.Range("J113").Value = 1905.2858
$ws.Range("K113").Value = 3718.3335
$ws.Range("L113").Value = 5715.857400000001
$ws.Range("M113").Value = -1548.3335
$ws.Range("N113").Value = -10055.8574
$ws.Range("H122").Value = 3822.8518
$ws.Range("I122").Value = 2254.923
$ws.Range("J122").Value = 5278.7856
$ws.Range("K122").Value = 6764.768999999999
$ws.Range("L122").Value = 15836.3568
$ws.Range("M122").Value = -4314.768999999999
$ws.Range("N122").Value = -20736.3568
$ws.Range("H124").Value = 37979.8
$ws.Range("J124").Value = 37979.8
$ws.Range("L124").Value = 37979.8
$ws.Range("N124").Value = -47799.8
$ws.Range("H126").Value = 1846.8
$ws.Range("I126").Value = 1325.2
$ws.Range("J126").Value = 2890
$ws.Range("K126").Value = 3975.6
$ws.Range("L126").Value = 8670
$ws.Range("M126").Value = -1505.6
$ws.Range("N126").Value = -13610
$ws.Range("H132").Value = 1438.0227
$ws.Range("I132").Value = 1414.9
$ws.Range("K132").Value = 4244.700000000001
$ws.Range("M132").Value = -1714.700000000001
